$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 379, shifting existing rows 379-438 down to 380-439
$ws.Rows("379:379").Insert()

# Populate the newly inserted row 379 with its data
$ws.Range("A379").Value = 4
$ws.Range("B379").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C379").Value = "Los Lagos"
$ws.Range("D379").Value = 44984
$ws.Range("E379").Value = 10
$ws.Range("F379").Value = 100114014
$ws.Range("G379").Value = "Betarraga"
$ws.Range("H379").Value = "Sin especificar"
$ws.Range("I379").Value = "Primera"
$ws.Range("J379").Value = 250
$ws.Range("K379").Value = 1000
$ws.Range("L379").Value = 1000
$ws.Range("M379").Value = 1000
$ws.Range("N379").Value = "$/paquete 5 unidades"
$ws.Range("O379").Value = "Provincia de Cautín"
$ws.Range("P379").Value = 200
$ws.Range("Q379").Value = 5
$ws.Range("R379").Value = "Hortaliza"
